$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1202
$ws1.Range("F5").Value = 629
$ws1.Range("F14").Value = 313
$ws1.Range("F19").Value = 12531
$ws1.Range("F20").Value = 12574

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 7

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1202
$ws4.Range("F6").Value = 629
$ws4.Range("F19").Value = 313
$ws4.Range("F23").Value = 7
$ws4.Range("F25").Value = 12531
$ws4.Range("F26").Value = 12574
